# Update sliding-window partial-model results (window 4) with new run output.
# Columns: A=Point, B=IPC RO (measured), C=IPC PO (predicted), D=DELTA (C-B), E=DELTA^2
# Row 52 = TOTAL (sum of DELTA, sum of DELTA^2), Row 53 = MSE (mean of DELTA^2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Cells.Item(2, 3).Value2 = 29.7437801361084   # C2
$ws.Cells.Item(2, 4).Value2 = 0.3737801361083939   # D2
$ws.Cells.Item(2, 5).Value2 = 0.1397115901492095   # E2

$ws.Cells.Item(3, 2).Value2 = 29.53999999999999   # B3
$ws.Cells.Item(3, 3).Value2 = 29.31538963317871   # C3
$ws.Cells.Item(3, 4).Value2 = -0.2246103668212811   # D3
$ws.Cells.Item(3, 5).Value2 = 0.05044981688359045   # E3

$ws.Cells.Item(4, 3).Value2 = 29.4737377166748   # C4
$ws.Cells.Item(4, 4).Value2 = -0.07626228332519247   # D4
$ws.Cells.Item(4, 5).Value2 = 0.005815935857971929   # E4

$ws.Cells.Item(5, 3).Value2 = 29.6823787689209   # C5
$ws.Cells.Item(5, 4).Value2 = -0.06762123107910156   # D5
$ws.Cells.Item(5, 5).Value2 = 0.004572630892653251   # E5

$ws.Cells.Item(6, 3).Value2 = 29.89393424987793   # C6
$ws.Cells.Item(6, 4).Value2 = 0.05393424987792628   # D6
$ws.Cells.Item(6, 5).Value2 = 0.00290890330989459   # E6

$ws.Cells.Item(7, 3).Value2 = 29.88800621032715   # C7
$ws.Cells.Item(7, 4).Value2 = 0.07800621032714616   # D7
$ws.Cells.Item(7, 5).Value2 = 0.006084968849602965   # E7

$ws.Cells.Item(8, 3).Value2 = 29.92011070251465   # C8
$ws.Cells.Item(8, 4).Value2 = 0.0001107025146467322   # D8
$ws.Cells.Item(8, 5).Value2 = 0.00000001225504674910996   # E8

$ws.Cells.Item(9, 3).Value2 = 30.01376152038574   # C9
$ws.Cells.Item(9, 4).Value2 = 0.03376152038573821   # D9
$ws.Cells.Item(9, 5).Value2 = 0.001139840258756617   # E9

$ws.Cells.Item(10, 2).Value2 = 30.03999999999999   # B10
$ws.Cells.Item(10, 3).Value2 = 30.16930961608887   # C10
$ws.Cells.Item(10, 4).Value2 = 0.1293096160888751   # D10
$ws.Cells.Item(10, 5).Value2 = 0.01672097681305228   # E10

$ws.Cells.Item(11, 2).Value2 = 30.21000000000001   # B11
$ws.Cells.Item(11, 3).Value2 = 30.22289848327637   # C11
$ws.Cells.Item(11, 4).Value2 = 0.01289848327635923   # D11
$ws.Cells.Item(11, 5).Value2 = 0.0001663708708305187   # E11

$ws.Cells.Item(12, 3).Value2 = 30.32418632507324   # C12
$ws.Cells.Item(12, 4).Value2 = 0.1041863250732433   # D12
$ws.Cells.Item(12, 5).Value2 = 0.01085479033226753   # E12

$ws.Cells.Item(13, 3).Value2 = 30.36539268493652   # C13
$ws.Cells.Item(13, 4).Value2 = -0.01460731506347202   # D13
$ws.Cells.Item(13, 5).Value2 = 0.0002133736533635364   # E13

$ws.Cells.Item(14, 3).Value2 = 30.54109001159668   # C14
$ws.Cells.Item(14, 4).Value2 = 0.101090011596682   # D14
$ws.Cells.Item(14, 5).Value2 = 0.01021919044461729   # E14

$ws.Cells.Item(15, 3).Value2 = 30.4333324432373   # C15
$ws.Cells.Item(15, 4).Value2 = -0.04666755676269929   # D15
$ws.Cells.Item(15, 5).Value2 = 0.00217786085419976   # E15

$ws.Cells.Item(16, 3).Value2 = 30.44302177429199   # C16
$ws.Cells.Item(16, 4).Value2 = -0.2469782257080055   # D16
$ws.Cells.Item(16, 5).Value2 = 0.06099824397387453   # E16

$ws.Cells.Item(17, 3).Value2 = 30.52885246276855   # C17
$ws.Cells.Item(17, 4).Value2 = -0.2211475372314453   # D17
$ws.Cells.Item(17, 5).Value2 = 0.04890623322353349   # E17

$ws.Cells.Item(18, 3).Value2 = 30.65643119812012   # C18
$ws.Cells.Item(18, 4).Value2 = -0.2835688018798805   # D18
$ws.Cells.Item(18, 5).Value2 = 0.08041126539959094   # E18

$ws.Cells.Item(19, 3).Value2 = 30.76356315612793   # C19
$ws.Cells.Item(19, 4).Value2 = -0.1864368438720732   # D19
$ws.Cells.Item(19, 5).Value2 = 0.03475869675297978   # E19

$ws.Cells.Item(20, 3).Value2 = 31.10555839538574   # C20
$ws.Cells.Item(20, 4).Value2 = 0.08555839538574617   # D20
$ws.Cells.Item(20, 5).Value2 = 0.007320239020983671   # E20

$ws.Cells.Item(21, 3).Value2 = 31.2354564666748   # C21
$ws.Cells.Item(21, 4).Value2 = 0.1154564666748001   # D21
$ws.Cells.Item(21, 5).Value2 = 0.01333019569702924   # E21

$ws.Cells.Item(22, 3).Value2 = 31.31316184997559   # C22
$ws.Cells.Item(22, 4).Value2 = 0.0331618499755848   # D22
$ws.Cells.Item(22, 5).Value2 = 0.001099708293803194   # E22

$ws.Cells.Item(23, 3).Value2 = 31.2291316986084   # C23
$ws.Cells.Item(23, 4).Value2 = -0.150868301391597   # D23
$ws.Cells.Item(23, 5).Value2 = 0.02276124436478575   # E23

$ws.Cells.Item(24, 3).Value2 = 31.40299797058105   # C24
$ws.Cells.Item(24, 4).Value2 = -0.1770020294189436   # D24
$ws.Cells.Item(24, 5).Value2 = 0.03132971841842458   # E24

$ws.Cells.Item(25, 2).Value2 = 31.65000000000001   # B25
$ws.Cells.Item(25, 3).Value2 = 31.89637565612793   # C25
$ws.Cells.Item(25, 4).Value2 = 0.246375656127924   # D25
$ws.Cells.Item(25, 5).Value2 = 0.06070096393246505   # E25

$ws.Cells.Item(26, 3).Value2 = 32.42288589477539   # C26
$ws.Cells.Item(26, 4).Value2 = 0.5428858947753952   # D26
$ws.Cells.Item(26, 5).Value2 = 0.2947250947460814   # E26

$ws.Cells.Item(27, 3).Value2 = 32.39987564086914   # C27
$ws.Cells.Item(27, 4).Value2 = 0.1198756408691395   # D27
$ws.Cells.Item(27, 5).Value2 = 0.0143701692737869   # E27

$ws.Cells.Item(28, 3).Value2 = 32.49711227416992   # C28
$ws.Cells.Item(28, 4).Value2 = 0.04711227416991903   # D28
$ws.Cells.Item(28, 5).Value2 = 0.00221956637746162   # E28

$ws.Cells.Item(29, 2).Value2 = 32.84999999999999   # B29
$ws.Cells.Item(29, 3).Value2 = 32.75205612182617   # C29
$ws.Cells.Item(29, 4).Value2 = -0.09794387817382244   # D29
$ws.Cells.Item(29, 5).Value2 = 0.009593003271728571   # E29

$ws.Cells.Item(30, 2).Value2 = 32.90000000000001   # B30
$ws.Cells.Item(30, 3).Value2 = 32.95458602905273   # C30
$ws.Cells.Item(30, 4).Value2 = 0.05458602905272869   # D30
$ws.Cells.Item(30, 5).Value2 = 0.002979634567745341   # E30

$ws.Cells.Item(31, 2).Value2 = 33.09999999999999   # B31
$ws.Cells.Item(31, 3).Value2 = 32.89218521118164   # C31
$ws.Cells.Item(31, 4).Value2 = -0.2078147888183537   # D31
$ws.Cells.Item(31, 5).Value2 = 0.04318698645161694   # E31

$ws.Cells.Item(32, 2).Value2 = 33.40000000000001   # B32
$ws.Cells.Item(32, 3).Value2 = 33.6392707824707   # C32
$ws.Cells.Item(32, 4).Value2 = 0.2392707824706974   # D32
$ws.Cells.Item(32, 5).Value2 = 0.05725050734413981   # E32

$ws.Cells.Item(33, 3).Value2 = 33.65039443969727   # C33
$ws.Cells.Item(33, 4).Value2 = -0.04960556030273722   # D33
$ws.Cells.Item(33, 5).Value2 = 0.002460711612948498   # E33

$ws.Cells.Item(34, 2).Value2 = 34.09999999999999   # B34
$ws.Cells.Item(34, 3).Value2 = 33.84643173217773   # C34
$ws.Cells.Item(34, 4).Value2 = -0.2535682678222599   # D34
$ws.Cells.Item(34, 5).Value2 = 0.06429686644638134   # E34

$ws.Cells.Item(35, 2).Value2 = 34.40000000000001   # B35
$ws.Cells.Item(35, 3).Value2 = 34.40230941772461   # C35
$ws.Cells.Item(35, 4).Value2 = 0.002309417724603691   # D35
$ws.Cells.Item(35, 5).Value2 = 0.000005333410226713688   # E35

$ws.Cells.Item(36, 2).Value2 = 34.90000000000001   # B36
$ws.Cells.Item(36, 3).Value2 = 35.04678726196289   # C36
$ws.Cells.Item(36, 4).Value2 = 0.1467872619628849   # D36
$ws.Cells.Item(36, 5).Value2 = 0.02154650027456061   # E36

$ws.Cells.Item(37, 3).Value2 = 35.71496200561523   # C37
$ws.Cells.Item(37, 4).Value2 = 0.4149620056152372   # D37
$ws.Cells.Item(37, 5).Value2 = 0.1721934661042202   # E37

$ws.Cells.Item(38, 3).Value2 = 36.00539016723633   # C38
$ws.Cells.Item(38, 4).Value2 = 0.3053901672363253   # D38
$ws.Cells.Item(38, 5).Value2 = 0.09326315424463072   # E38

$ws.Cells.Item(39, 3).Value2 = 36.00448989868164   # C39
$ws.Cells.Item(39, 4).Value2 = -0.2955101013183565   # D39
$ws.Cells.Item(39, 5).Value2 = 0.08732621998118534   # E39

$ws.Cells.Item(40, 3).Value2 = 36.59141159057617   # C40
$ws.Cells.Item(40, 4).Value2 = -0.2085884094238253   # D40
$ws.Cells.Item(40, 5).Value2 = 0.04350912454596136   # E40

$ws.Cells.Item(41, 3).Value2 = 37.24056625366211   # C41
$ws.Cells.Item(41, 4).Value2 = -0.05943374633788778   # D41
$ws.Cells.Item(41, 5).Value2 = 0.003532370203756389   # E41

$ws.Cells.Item(42, 2).Value2 = 37.90000000000001   # B42
$ws.Cells.Item(42, 3).Value2 = 37.99029922485352   # C42
$ws.Cells.Item(42, 4).Value2 = 0.09029922485350994   # D42
$ws.Cells.Item(42, 5).Value2 = 0.008153950009144748   # E42

$ws.Cells.Item(43, 3).Value2 = 38.42107772827148   # C43
$ws.Cells.Item(43, 4).Value2 = -0.07892227172851562   # D43
$ws.Cells.Item(43, 5).Value2 = 0.006228724974789657   # E43

$ws.Cells.Item(44, 2).Value2 = 38.90000000000001   # B44
$ws.Cells.Item(44, 3).Value2 = 39.00311660766602   # C44
$ws.Cells.Item(44, 4).Value2 = 0.1031166076660099   # D44
$ws.Cells.Item(44, 5).Value2 = 0.01063303477654582   # E44

$ws.Cells.Item(45, 2).Value2 = 39.40000000000001   # B45
$ws.Cells.Item(45, 3).Value2 = 39.53226852416992   # C45
$ws.Cells.Item(45, 4).Value2 = 0.1322685241699162   # D45
$ws.Cells.Item(45, 5).Value2 = 0.0174949624860877   # E45

$ws.Cells.Item(46, 2).Value2 = 39.90000000000001   # B46
$ws.Cells.Item(46, 3).Value2 = 39.5301399230957   # C46
$ws.Cells.Item(46, 4).Value2 = -0.3698600769043026   # D46
$ws.Cells.Item(46, 5).Value2 = 0.1367964764876566   # E46

$ws.Cells.Item(47, 2).Value2 = 40.09999999999999   # B47
$ws.Cells.Item(47, 3).Value2 = 39.9373664855957   # C47
$ws.Cells.Item(47, 4).Value2 = -0.1626335144042912   # D47
$ws.Cells.Item(47, 5).Value2 = 0.02644966000749079   # E47

$ws.Cells.Item(48, 2).Value2 = 40.59999999999999   # B48
$ws.Cells.Item(48, 3).Value2 = 40.45168685913086   # C48
$ws.Cells.Item(48, 4).Value2 = -0.1483131408691349   # D48
$ws.Cells.Item(48, 5).Value2 = 0.02199678775446786   # E48

$ws.Cells.Item(49, 2).Value2 = 40.90000000000001   # B49
$ws.Cells.Item(49, 3).Value2 = 40.73793411254883   # C49
$ws.Cells.Item(49, 4).Value2 = -0.1620658874511776   # D49
$ws.Cells.Item(49, 5).Value2 = 0.02626535187533775   # E49

$ws.Cells.Item(50, 2).Value2 = 41.20000000000001   # B50
$ws.Cells.Item(50, 3).Value2 = 41.31341552734375   # C50
$ws.Cells.Item(50, 4).Value2 = 0.1134155273437401   # D50
$ws.Cells.Item(50, 5).Value2 = 0.01286308184265865   # E50

$ws.Cells.Item(51, 3).Value2 = 41.75338745117188   # C51
$ws.Cells.Item(51, 4).Value2 = 0.253387451171875   # D51
$ws.Cells.Item(51, 5).Value2 = 0.06420520041137934   # E51

$ws.Cells.Item(52, 3).Value2 = 0.1432662963866917   # C52
$ws.Cells.Item(52, 5).Value2 = 1.856198709984518   # E52

$ws.Cells.Item(53, 5).Value2 = 0.03712397419969036   # E53
